# Update "想去人数" (F column) and one "最低票价" (G19) values
# on both the "展览" and "全部类型" worksheets, which hold identical data.

$wb = $excel.ActiveWorkbook

# Cell -> new value map (same for both affected sheets)
$updates = @{
    "F4"  = 255
    "F5"  = 52
    "F6"  = 177
    "F8"  = 38
    "F11" = 40
    "F13" = 93
    "F14" = 1426
    "F16" = 488
    "F17" = 439
    "G19" = 35
    "F20" = 35
    "F21" = 39
    "F22" = 1352
    "F23" = 3308
    "F25" = 54
    "F27" = 1071
    "F28" = 73
    "F29" = 1661
    "F32" = 47
    "F33" = 274
    "F34" = 399
    "F35" = 451
    "F36" = 626
    "F38" = 24
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
